$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.905.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.879.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4613"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3868"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07857"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9859"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.77"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.891.02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.989"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.647"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06978"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009951"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.908.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.246"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.104"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.990"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "117.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.911"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09364"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9008"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.260"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.319"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.256"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.179"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05744"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02075"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("E39").Value = "  -5.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5642"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1766"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.712"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.273"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5343"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07041"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.04%  "
$ws.Range("E47").Value = "  -1.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.553"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "112.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.066"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.76%  "
